$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.104.97"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "2.166.00"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'252.99"
$ws.Range("E5").Value = "  +6.17%  "

$ws.Range("D6").Value = "'0.607"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'73.18"
$ws.Range("E7").Value = "  +1.49%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").Value = "'39.80"
$ws.Range("E10").Value = "  -0.28%  "

$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = "  +0.51%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.101"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.75"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").Value = "2.493.06"
$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("D15").Value = "'14.21"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "2.173.71"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("D17").Value = "'0.769"
$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("D18").Value = "42.024.66"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("E19").Value = "  -2.13%  "

$ws.Range("D20").Value = "'70.73"
$ws.Range("E20").Value = "  +1.22%  "

$ws.Range("D21").Value = "'5.84"
$ws.Range("E21").Value = "  +1.18%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'9.58"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'2.17"
$ws.Range("E23").Value = "  +6.43%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'226.19"
$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").Value = "'10.52"
$ws.Range("E26").Value = "  -1.85%  "

$ws.Range("E27").Value = "  +1.59%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  +2.76%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("D30").Value = "'36.90"
$ws.Range("E30").Value = "  +12.01%  "

$ws.Range("D31").Value = "'168.72"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").Value = "'19.90"
$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("D33").Value = "'0.0800"
$ws.Range("E33").Value = "  +4.13%  "

$ws.Range("D34").Value = "'5.11"
$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("D36").Value = "'0.106"
$ws.Range("E36").Value = "  +1.90%  "

$ws.Range("D37").Value = "'4.27"
$ws.Range("E37").Value = "  -0.91%  "

$ws.Range("D38").Value = "'0.0326"
$ws.Range("E38").Value = "  +7.76%  "

$ws.Range("D39").Value = "'12.01"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("D40").Value = "'2.06"
$ws.Range("E40").Value = "  -2.10%  "

$ws.Range("D41").Value = "'0.196"
$ws.Range("E41").Value = "  +3.90%  "

$ws.Range("D42").Value = "'5.15"
$ws.Range("E42").Value = "  -3.44%  "

$ws.Range("D43").Value = "'58.68"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("D44").Value = "'100.86"
$ws.Range("E44").Value = "  +3.99%  "

$ws.Range("D45").Value = "'0.470"
$ws.Range("E45").Value = "  +20.81%  "

$ws.Range("D46").Value = "'8.25"
$ws.Range("E46").Value = "  -1.90%  "

$ws.Range("D47").Value = "'0.0965"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("D48").Value = "'2.39"
$ws.Range("E48").Value = "  +9.61%  "

$ws.Range("D49").Value = "'1.09"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").Value = "'1.12"
$ws.Range("E50").Value = "  +1.31%  "

$ws.Range("E51").Value = "  +1.18%  "

Write-Output "Updated cryptos list"